# Auto-generated script to apply cryptos.xlsx diff via Excel COM interop
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$apostrophe = [char]39

$ws.Range('D2').Value = '69.843.21'
$ws.Range('E2').Value = '  -1.71%  '
$ws.Range('D3').Value = '3.558.81'
$ws.Range('E3').Value = '  -2.81%  '
$ws.Range('D4').Value = "$apostrophe" + '1.00'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = "$apostrophe" + '574.77'
$ws.Range('E5').Value = '  -3.61%  '
$ws.Range('D6').Value = "$apostrophe" + '185.32'
$ws.Range('E6').Value = '  -4.82%  '
$ws.Range('D7').Value = '3.551.75'
$ws.Range('E7').Value = '  -2.80%  '
$ws.Range('D8').Value = "$apostrophe" + '0.618'
$ws.Range('E8').Value = '  -4.67%  '
$ws.Range('E9').Value = '  +0.10%  '
$ws.Range('D10').Value = "$apostrophe" + '0.183'
$ws.Range('E10').Value = '  -1.02%  '
$ws.Range('D11').Value = "$apostrophe" + '0.646'
$ws.Range('E11').Value = '  -4.31%  '
$ws.Range('D12').Value = "$apostrophe" + '54.69'
$ws.Range('E12').Value = '  -6.44%  '
$ws.Range('E13').Value = '  +1.46%  '
$ws.Range('E14').Value = '  -4.95%  '
$ws.Range('D15').Value = '4.129.11'
$ws.Range('E15').Value = '  -2.66%  '
$ws.Range('D16').Value = "$apostrophe" + '19.50'
$ws.Range('E16').Value = '  -3.15%  '
$ws.Range('D17').Value = '3.554.61'
$ws.Range('E17').Value = '  -2.66%  '
$ws.Range('D18').Value = '69.776.65'
$ws.Range('E18').Value = '  -1.67%  '
$ws.Range('D19').Value = "$apostrophe" + '12.51'
$ws.Range('E19').Value = '  -2.71%  '
$ws.Range('E20').Value = '  -0.96%  '
$ws.Range('E21').Value = '  -4.08%  '
$ws.Range('D22').Value = "$apostrophe" + '490.86'
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('D23').Value = "$apostrophe" + '19.22'
$ws.Range('E23').Value = '  +0.47%  '
$ws.Range('D24').Value = "$apostrophe" + '4.87'
$ws.Range('E24').Value = '  -8.03%  '
$ws.Range('E25').Value = '  -3.61%  '
$ws.Range('D26').Value = "$apostrophe" + '94.91'
$ws.Range('E26').Value = '  +3.61%  '
$ws.Range('D27').Value = "$apostrophe" + '11.44'
$ws.Range('E27').Value = '  -0.79%  '
$ws.Range('E28').Value = '  -7.36%  '
$ws.Range('D29').Value = "$apostrophe" + '9.25'
$ws.Range('E29').Value = '  -4.09%  '
$ws.Range('D30').Value = "$apostrophe" + '31.52'
$ws.Range('E30').Value = '  -4.33%  '
$ws.Range('D31').Value = "$apostrophe" + '7.45'
$ws.Range('E31').Value = '  -4.33%  '
$ws.Range('E32').Value = '  -0.20%  '
$ws.Range('E33').Value = '  -2.75%  '
$ws.Range('D34').Value = "$apostrophe" + '0.115'
$ws.Range('E34').Value = '  -6.41%  '
$ws.Range('D35').Value = "$apostrophe" + '564.74'
$ws.Range('E35').Value = '  -10.38%  '
$ws.Range('E36').Value = '  +11.09%  '
$ws.Range('D37').Value = "$apostrophe" + '38.61'
$ws.Range('E37').Value = '  -4.27%  '
$ws.Range('E38').Value = '  -0.05%  '
$ws.Range('E39').Value = '  -5.24%  '
$ws.Range('E40').Value = '  -5.50%  '
$ws.Range('D41').Value = "$apostrophe" + '3.50'
$ws.Range('E41').Value = '  -2.67%  '
$ws.Range('D42').Value = "$apostrophe" + '3.14'
$ws.Range('E42').Value = '  +2.74%  '
$ws.Range('E44').Value = '  -6.21%  '
$ws.Range('D45').Value = '3.209.95'
$ws.Range('E45').Value = '  -3.13%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').Value = "$apostrophe" + '0.0438'
$ws.Range('E46').Value = '  -4.46%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').Value = "$apostrophe" + '3.44'
$ws.Range('E47').Value = '  +3.81%  '
$ws.Range('E48').Value = '  -0.21%  '
$ws.Range('E49').Value = '  -3.08%  '
$ws.Range('D50').Value = "$apostrophe" + '0.999'
$ws.Range('E50').Value = '  +0.13%  '
$ws.Range('D51').Value = "$apostrophe" + '3.12'
$ws.Range('E51').Value = '  -5.32%  '

Write-Host "Applied" 86 "cell updates"
